$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '64.118.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -1.42%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.521.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.03%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.05%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '585.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -1.17%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '133.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.48%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '3.521.62'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.03%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.80%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -0.15%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -0.51%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -1.76%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '4.127.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +0.15%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -0.34%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +1.35%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = "'" + 'WrappedEther'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = "'" + '3.533.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +0.22%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = "'" + 'ShibaInu'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'" + 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'" + '0.0000179'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -1.57%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '64.176.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -1.26%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -3.47%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -2.57%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -1.13%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '385.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -1.69%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -1.01%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '3.667.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.05%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '74.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -1.06%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -0.01%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '0.0000115'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +2.54%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '1.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -2.41%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '7.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -3.45%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.01%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '8.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +1.01%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -2.25%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '3.534.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +0.15%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +0.00%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '23.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -2.19%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +1.32%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '5.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +0.84%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -0.58%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '6.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -0.72%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '160.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -4.43%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.0789'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -2.63%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -1.16%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '26.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +1.95%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +0.12%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '41.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -3.12%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -4.69%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '4.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -0.42%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -3.00%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '2.473.78'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +1.86%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -1.37%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.908'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -0.32%  '
$ws.Range('E51').Style = 'Normal'
